# Update Front-End for Products page
# - "III. Project Issues" table: the item that used to be row 20
#   ("Discussion, analysis of ERD diagram." / "Pending") is moved up into
#   row 19 (replacing "Learn API"), and row 20 is cleared out, with its
#   B cell re-filled white (no border) instead of the bordered table style.
# - The frozen-pane view is scrolled down and the active cell is reset to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the "Learn API" row content with the text that used to live in
# row 20, then blank out row 20's B/C cells.
$ws.Range("B19").Value = "Discussion, analysis of ERD diagram."
$ws.Range("B20").ClearContents()
$ws.Range("C20").ClearContents()

# Give the now-empty B20 cell a plain white fill with no border (instead of
# the bordered "table" style it used to share with the rest of the grid).
$ws.Range("B20").Borders.LineStyle = -4142
$ws.Range("B20").Interior.Color = 16777215

# Restore the frozen pane (keep the freeze after row 3) and scroll the
# view down so row 8 is the first visible row under the frozen header,
# then set the active selection to B1.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A4").Select()
$win.FreezePanes = $true
$win.ScrollRow = 8
$win.ScrollColumn = 1

$ws.Range("B1").Select()
